$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.782.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.10%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.408.74"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.15%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "568.68"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.50%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.69"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.59%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.414.24"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.21%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.579"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +9.28%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.34"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.98%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.70%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.443"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.28%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.989.90"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.02%  "

$ws.Range("E14").Value = "  -3.10%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000191"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.73%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.80%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.686.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.89%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.441.59"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.14%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.39"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.52%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.09%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "383.02"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.84%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "8.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -3.96%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "73.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.999"
$ws.Range("D24").Style = "Normal"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.538"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.51%  "

$ws.Range("E26").Value = "  +22.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.57%  "

$ws.Range("E28").Value = "  -0.36%  "

$ws.Range("E29").Value = "  +0.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +9.62%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.40"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +6.09%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.00"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.60%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "23.41"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.56%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.43"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.29%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.92"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +2.71%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "160.80"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.86%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.46"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.84%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0764"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.87%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.85"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.49%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.894.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.36%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.81"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.91%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0318"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.21%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.17"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +3.16%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.37%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.757"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.32%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "23.27"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +5.62%  "

$ws.Range("E48").Value = "  +2.26%  "

$ws.Range("B49").Value = "Stellar"
$ws.Range("C49").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.108"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.24%  "

$ws.Range("B50").Value = "dogwifhat"
$ws.Range("C50").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.17"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +19.60%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.48"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.42%  "

Write-Host "Updated cryptos list"
